$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Name = "NPC_01"
$ws.Activate()
$ws.Range("C1").Select()
